# Initial version of version 2.0.0.
# Adds IBD0/IBD1/IBD2 columns (with worked fractions) and two extra
# relatedness categories to the "Relatedness" sheet; tidies up the
# column widths on "Relatedness" and "PheWAS"; and moves the active
# selection from "PheWAS" to "Relatedness".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Relatedness")
$wsPheWAS = $wb.Worksheets.Item("PheWAS")

# --- Header row: new IBD0 / IBD1 / IBD2 columns -----------------------
$ws.Range("C1").Value = "IBD0"
$ws.Range("D1").Value = "IBD1"
$ws.Range("E1").Value = "IBD2"

# --- Row 2: Monozygotic twins ----------------------------------------
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0

# --- Row 3: Parents/child ---------------------------------------------
$ws.Range("C3").Value = 0.25
$ws.Range("D3").Value = 0.5
$ws.Range("E3").Value = 0.25
$ws.Range("D3").NumberFormat = "0.00"

# --- Row 4: Sibling -----------------------------------------------------
$ws.Range("C4").Value = 0.25
$ws.Range("D4").Value = 0.5
$ws.Range("E4").Value = 0.25
$ws.Range("D4").NumberFormat = "0.00"

# --- Row 5: Fraternal twins --------------------------------------------
$ws.Range("C5").Value = 0.25
$ws.Range("D5").Value = 0.5
$ws.Range("E5").Value = 0.25
$ws.Range("D5").NumberFormat = "0.00"

# --- Row 6: Grandparent/grandchild -------------------------------------
$ws.Range("C6").Value = 0.5
$ws.Range("D6").Value = 0.5
$ws.Range("E6").Value = 0
$ws.Range("C6").NumberFormat = "0.00"
$ws.Range("D6").NumberFormat = "0.00"

# --- Row 7: Aunt/Uncle/Niece/Nephew ------------------------------------
$ws.Range("C7").Value = 0.5
$ws.Range("D7").Value = 0.5
$ws.Range("E7").Value = 0
$ws.Range("C7").NumberFormat = "0.00"
$ws.Range("D7").NumberFormat = "0.00"

# --- Row 8: Half-sibling -------------------------------------------------
$ws.Range("C8").Value = 0.5
$ws.Range("D8").Value = 0.5
$ws.Range("E8").Value = 0
$ws.Range("C8").NumberFormat = "0.00"
$ws.Range("D8").NumberFormat = "0.00"

# --- Row 9: First-cousin (formulas) --------------------------------------
$ws.Range("C9").Formula = "=C8+D9"
$ws.Range("D9").Formula = "=D8/2"
$ws.Range("E9").Value = 0

# --- Row 10: Half first-cousin -------------------------------------------
$ws.Range("C10").Formula = "=C9+D10"
$ws.Range("D10").Formula = "=D9/2"
$ws.Range("E10").Value = 0

# --- Row 11: First-cousin once removed -----------------------------------
$ws.Range("C11").Formula = "=C10"
$ws.Range("D11").Formula = "=D10"
$ws.Range("E11").Value = 0

# --- Row 12: Second-cousin -------------------------------------------------
$ws.Range("C12").Formula = "=0.75+0.125+D12"
$ws.Range("D12").Formula = "=D11/2"
$ws.Range("E12").Value = 0

# --- Row 13: Second-cousin once removed ------------------------------------
$ws.Range("C13").Formula = "=C12+D13"
$ws.Range("D13").Formula = "=D12/2"
$ws.Range("E13").Value = 0

# --- New rows 14 & 15: Distantly related / Unrelated -----------------------
# (written in this order so shared-string insertion order matches)
$ws.Range("A15").Value = "Unrelated (includes relationships beyond the third degree)"
$ws.Range("A14").Value = "Distantly related"
$ws.Range("B15").Value = "<1.56%"
$ws.Range("B14").Value = "<1.56%"
$ws.Range("C14").Value = "varies"
$ws.Range("D14").Value = "varies"
$ws.Range("E14").Value = 0
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0

# --- Column widths -----------------------------------------------------
$ws.Columns.Item("A").ColumnWidth = 25.1666666
$ws.Columns.Item("B").ColumnWidth = 13.1666666
$wsPheWAS.Columns.Item("A").ColumnWidth = 36.3333333

# --- Selection / active sheet -------------------------------------------
$ws.Activate()
$ws.Range("F10").Select()
